$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: optimization_parameters
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1: drop the extra "value" header cells in C1:F1 (only A1/B1 remain).
$ws.Range("C1:F1").ClearContents()

# Row 8 "Model" -> "production_function" (keep the Sigmoid value in B8).
$ws.Range("A8").Value = "production_function"

# Insert a brand-new row above the old row 9 (estimate_params) for the new
# "L_curve" parameter.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("B9").NumberFormat = "0.000"

# Remove the old "Deletion" row entirely (now shifted to row 17 after the
# insert above).
$ws.Rows.Item(17).Delete()

# Sheet7 becomes the active/selected sheet, with the whole last row selected.
$ws.Rows.Item(17).Select()

# ---------------------------------------------------------------------------
# Sheet: network_weights loses the tab-selection to optimization_parameters.
# ---------------------------------------------------------------------------
$wsw = $wb.Worksheets.Item("network_weights")
$wsw.Range("A1:E5").Select()

$ws.Activate()

$wb.Save()
